$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2: base EXPON.DIST formula (mean = 1/E1 -> lambda)
$ws.Range("B2").Formula = '=EXPON.DIST(A2, 1/$E$1,TRUE)'

# Fill the formula down the column in the same chunks the original
# tutorial used, producing the same shared-formula groupings.
$ws.Range("B3:B66").Formula = '=EXPON.DIST(A3, 1/$E$1,TRUE)'
$ws.Range("B67:B130").Formula = '=EXPON.DIST(A67, 1/$E$1,TRUE)'
$ws.Range("B131:B151").Formula = '=EXPON.DIST(A131, 1/$E$1,TRUE)'

# E3: difference between the probability at 15 seconds and at 10 seconds
$ws.Range("E3").Formula = '=B16-B11'

# Leave the selection on E4, matching the saved workbook state
$ws.Range("E4").Select()
